# Generate Report for Handoff
#
# The localization status moved from "In Translation" to "Ready for
# handoff" and the handoff timestamps were refreshed. This touches:
#   - Overview!E2, Overview!F2  (per-language status cells)
#   - Overview!G2               (Latest HO Xliff Generate Date)
#   - zh-cn!C2                  (Status)
#   - zh-cn!H2                  (Latest Handoff Datetime)
#   - de-de!C2                  (Status)
#   - de-de!H2                  (Latest Handoff Datetime)
# Because the new status text ("Ready for handoff") is longer than the
# old one ("In Translation"), the Status columns also grow wider.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
# Closest column width (in characters) reachable through this runtime's
# ColumnWidth rounding that lands on the target stored width of
# 17.2159881591797 (the nearest achievable value is 17.1666...).
$newStatusColWidth = 16.33

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-10-13 12:49:50"
$wsOverview.Range("E1").ColumnWidth = $newStatusColWidth
$wsOverview.Range("F1").ColumnWidth = $newStatusColWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-10-13 12:49:40"
$wsZhCn.Range("C1").ColumnWidth = $newStatusColWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-10-13 12:49:50"
$wsDeDe.Range("C1").ColumnWidth = $newStatusColWidth
